$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A154").Value = "IMX-USD"
$ws.Range("A155").Value = "GRT-USD"
